$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Remove the "Uhrzeit" (time) column - this shifts Schul/Uni..Teilnehmer one column left
$ws.Columns.Item(2).Delete()

# Apply Text format to the (now shifted) Schul/Uni, Adresse, Stadt, Bundesland, PLZ columns
$ws.Columns.Item("B:F").NumberFormat = "@"

# Apply integer Number format to the Tische / Teilnehmer columns
$ws.Columns.Item("G:H").NumberFormat = "0"

# Selection left behind by the formatting action
$ws.Range("B1:E1048576").Select()
